$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.081.19"
$ws.Range("E2").Value = "  +8.64%  "

$ws.Range("D3").Value = "3.453.11"
$ws.Range("E3").Value = "  +5.94%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "414.58"
$ws.Range("E5").Value = "  +4.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "125.02"
$ws.Range("E6").Value = "  +14.95%  "

$ws.Range("D7").Value = "3.449.77"
$ws.Range("E7").Value = "  +5.98%  "

$ws.Range("E8").Value = "  +2.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.661"
$ws.Range("E10").Value = "  +6.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.128"
$ws.Range("E11").Value = "  +34.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.53"
$ws.Range("E12").Value = "  +5.48%  "

$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "3.994.07"
$ws.Range("E14").Value = "  +5.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.53"
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.77"
$ws.Range("E16").Value = "  +4.15%  "

$ws.Range("D17").Value = "3.434.40"
$ws.Range("E17").Value = "  +5.35%  "

$ws.Range("D18").Value = "62.014.55"
$ws.Range("E18").Value = "  +8.77%  "

$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000133"
$ws.Range("E21").Value = "  +22.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.35"
$ws.Range("E22").Value = "  +0.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.10"
$ws.Range("E23").Value = "  +10.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "314.23"
$ws.Range("E24").Value = "  +6.87%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.03"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "31.00"
$ws.Range("E27").Value = "  +10.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.84"
$ws.Range("E28").Value = "  +5.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.89"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  +4.25%  "

$ws.Range("E33").Value = "  +3.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("E34").Value = "  +20.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "42.22"
$ws.Range("E35").Value = "  +5.08%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.32"
$ws.Range("E38").Value = "  +1.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.02"
$ws.Range("E42").Value = "  +7.98%  "

$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.44"
$ws.Range("E44").Value = "  -1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.27"
$ws.Range("E45").Value = "  +2.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.286"
$ws.Range("E46").Value = "  +0.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.90"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.32"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("D50").Value = "2.211.34"
$ws.Range("E50").Value = "  +2.93%  "

$ws.Range("D51").Value = "3.790.34"
$ws.Range("E51").Value = "  +5.69%  "
